$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values that changed
$ws.Range("F2").Value = 0.5200242194585092
$ws.Range("E3").Value = 0.5608680360343469

# Delete rows 4 to 6 (removing the extra data rows entirely)
$ws.Range("A4:I6").EntireRow.Delete()
